$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 85523
$ws.Range("B2").Value = "Ana Liz Lima"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45101
$ws.Range("G2").Value = 3187.25

# Row 3
$ws.Range("A3").Value = 35015
$ws.Range("B3").Value = "Sara Mendonça"
$ws.Range("C3").Value = "TI"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 45102
$ws.Range("G3").Value = 6666.56

# Row 4
$ws.Range("A4").Value = 28014
$ws.Range("B4").Value = "Ravi Lucca Costa"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 4969.03

# Row 5
$ws.Range("A5").Value = 36475
$ws.Range("B5").Value = "Heitor da Luz"
$ws.Range("C5").Value = "Financeiro"
$ws.Range("D5").Value = "Doenca"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 45089
$ws.Range("G5").Value = 7505.13

# Row 6
$ws.Range("A6").Value = 61840
$ws.Range("B6").Value = "Vinicius Oliveira"
$ws.Range("C6").Value = "Juridico"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 9727.65

# Row 7
$ws.Range("A7").Value = 80193
$ws.Range("B7").Value = "Maysa Leão"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45106
$ws.Range("G7").Value = 8282.309999999999

# Row 8
$ws.Range("A8").Value = 40673
$ws.Range("B8").Value = "Luna Carvalho"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Outros"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45086
$ws.Range("G8").Value = 3385.38

# Row 9
$ws.Range("A9").Value = 19875
$ws.Range("B9").Value = "Marcos Vinicius Gonçalves"
$ws.Range("C9").Value = "Juridico"
$ws.Range("D9").Value = "Consulta medica"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 45078
$ws.Range("G9").Value = 5643.04

# Row 10
$ws.Range("A10").Value = 29406
$ws.Range("B10").Value = "Luiz Miguel Castro"
$ws.Range("C10").Value = "Operacoes"
$ws.Range("D10").Value = "Viagem de negocios"
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45091
$ws.Range("G10").Value = 9160.450000000001

# Row 11
$ws.Range("A11").Value = 93862
$ws.Range("B11").Value = "Evelyn Vieira"
$ws.Range("C11").Value = "Operacoes"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45098
$ws.Range("G11").Value = 9628.280000000001
